$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" row total correct count
$ws.Range("B11").Value = 5

# Update the "Total" row correct count and the corresponding fraction text
$ws.Range("B12").Value = 90
$ws.Range("E12").Value = "90/140"
